{"js": "// Update codes to integrate sensitivity analyses\n//\n// 1) Add a new \"Abstract Title\" paragraph style (based on Normal, followed\n//    by Abstract), centered/bold/colored small caption-like heading used\n//    right before the Abstract paragraph.\n// 2) Tighten the space-before on the existing \"Abstract\" style (300 -> 100)\n//    since the new Abstract Title style now carries the spacing above it.\n// 3) Add a new \"Footnote Block Text\" paragraph style (based on Footnote\n//    Text) for indented block quotations inside footnotes.\n\n// --- 1) AbstractTitle --------------------------------------------------\ncontext.document.addStyle(\"Abstract Title\", Word.StyleType.paragraph);\nawait context.sync();\n\nconst abstractTitle = context.document.getStyles().getByName(\"Abstract Title\");\nabstractTitle.baseStyle = \"Normal\";\nabstractTitle.nextParagraphStyle = \"Abstract\";\nabstractTitle.quickStyle = true;\nabstractTitle.paragraphFormat.keepWithNext = true;\nabstractTitle.paragraphFormat.keepTogether = true;\nabstractTitle.paragraphFormat.alignment = Word.Alignment.centered;\nabstractTitle.paragraphFormat.spaceBefore = 15;\nabstractTitle.paragraphFormat.spaceAfter = 0;\nabstractTitle.font.size = 10;\nabstractTitle.font.sizeBidirectional = 10;\nabstractTitle.font.bold = true;\nabstractTitle.font.color = \"#345A8A\";\n\n// --- 2) Abstract: before-spacing 300 -> 100 -----------------------------\nconst abstract = context.document.getStyles().getByName(\"Abstract\");\nabstract.paragraphFormat.spaceBefore = 5;\n\nawait context.sync();\n\n// --- 3) FootnoteBlockText -----------------------------------------------\ncontext.document.addStyle(\"Footnote Block Text\", Word.StyleType.paragraph);\nawait context.sync();\n\nconst footnoteBlockText = context.document.getStyles().getByName(\"Footnote Block Text\");\nfootnoteBlockText.baseStyle = \"Footnote Text\";\nfootnoteBlockText.nextParagraphStyle = \"Footnote Text\";\nfootnoteBlockText.priority = 9;\nfootnoteBlockText.unhideWhenUsed = true;\nfootnoteBlockText.quickStyle = true;\nfootnoteBlockText.paragraphFormat.spaceBefore = 5;\nfootnoteBlockText.paragraphFormat.spaceAfter = 5;\nfootnoteBlockText.paragraphFormat.firstLineIndent = 0;\nfootnoteBlockText.paragraphFormat.leftIndent = 24;\nfootnoteBlockText.paragraphFormat.rightIndent = 24;\n\nawait context.sync();\n", "ps1": "# Update codes to integrate sensitivity analyses\n#\n# 1) Add a new \"Abstract Title\" paragraph style (based on Normal, followed\n#    by Abstract), centered/bold/colored small caption-like heading used\n#    right before the Abstract paragraph.\n# 2) Tighten the space-before on the existing \"Abstract\" style (300 -> 100)\n#    since the new Abstract Title style now carries the spacing above it.\n# 3) Add a new \"Footnote Block Text\" paragraph style (based on Footnote\n#    Text) for indented block quotations inside footnotes.\n\n$d = $word.ActiveDocument\n\n# --- 1) AbstractTitle ----------------------------------------------------\n$abstractTitle = $d.Styles.Add(\"Abstract Title\", 1)\n$abstractTitle.BaseStyle = $d.Styles(\"Normal\")\n$abstractTitle.NextParagraphStyle = $d.Styles(\"Abstract\")\n$abstractTitle.QuickStyle = $true\n$abstractTitle.ParagraphFormat.KeepWithNext = $true\n$abstractTitle.ParagraphFormat.KeepTogether = $true\n$abstractTitle.ParagraphFormat.Alignment = 1\n$abstractTitle.ParagraphFormat.SpaceBefore = 15\n$abstractTitle.ParagraphFormat.SpaceAfter = 0\n$abstractTitle.Font.Size = 10\n$abstractTitle.Font.SizeBi = 10\n$abstractTitle.Font.Bold = $true\n$abstractTitle.Font.Color = 9067060\n\n# --- 2) Abstract: before-spacing 300 -> 100 ------------------------------\n$abstract = $d.Styles(\"Abstract\")\n$abstract.ParagraphFormat.SpaceBefore = 5\n\n# --- 3) FootnoteBlockText -------------------------------------------------\n$footnoteBlockText = $d.Styles.Add(\"Footnote Block Text\", 1)\n$footnoteBlockText.BaseStyle = $d.Styles(\"Footnote Text\")\n$footnoteBlockText.NextParagraphStyle = $d.Styles(\"Footnote Text\")\n$footnoteBlockText.Priority = 9\n$footnoteBlockText.UnhideWhenUsed = $true\n$footnoteBlockText.QuickStyle = $true\n$footnoteBlockText.ParagraphFormat.SpaceBefore = 5\n$footnoteBlockText.ParagraphFormat.SpaceAfter = 5\n$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0\n$footnoteBlockText.ParagraphFormat.LeftIndent = 24\n$footnoteBlockText.ParagraphFormat.RightIndent = 24\n"}
